$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Title ----
$ws.Range("A1").Value = "龙场营镇者把小学学生营养餐食品原材料采购台账"

# Insert 4 new rows before the footer (old row 7), pushing it down to row 11,
# carrying the formatting of the row above (row 6) into the newly inserted
# rows 7-10, matching Excel's normal row-insert formatting behavior.
$ws.Rows("7:10").Insert()

# The 数量 (E) / 单价 (G) / 金额 (H) columns hold numeric-looking values that
# must be stored as text (as in the source export), so force Text format on
# them before writing so Excel doesn't silently coerce them to numbers.
$ws.Range("E4:E10").NumberFormat = "@"
$ws.Range("G4:G10").NumberFormat = "@"
$ws.Range("H4:H10").NumberFormat = "@"

# ---- Row 4 ----
$ws.Range("A4").Value = "Fri Jul 06 2018 00:51:15 GMT+0800 (China Standard Time)"
$ws.Range("B4").Value = "瘦肉"
$ws.Range("C4").Value = "Fri Jul 06 2018 00:51:15 GMT+0800 (China Standard Time)"
$ws.Range("D4").Value = "Fri Jul 06 2018 00:51:15 GMT+0800 (China Standard Time)"
$ws.Range("E4").Value = "30"
$ws.Range("F4").Value = "斤"
$ws.Range("G4").Value = "15"
$ws.Range("H4").Value = "450"
$ws.Range("I4").Value = "YU"
$ws.Range("J4").Value = "YU"
$ws.Range("K4").Value = "YU"

# ---- Row 5 ----
$ws.Range("A5").Value = "Sun Jul 08 2018 23:05:33 GMT+0800 (China Standard Time)"
$ws.Range("B5").Value = "瘦肉"
$ws.Range("E5").Value = "1"
$ws.Range("F5").Value = "斤"
$ws.Range("G5").Value = "15"
$ws.Range("H5").Value = "15"
$ws.Range("I5").Value = "采购人"
$ws.Range("J5").Value = "收验货人"
$ws.Range("K5").Value = "供货人"

# ---- Row 6 ----
$ws.Range("A6").Value = "Mon Jul 09 2018 23:43:52 GMT+0800 (China Standard Time)"
$ws.Range("B6").Value = "瘦肉"
$ws.Range("E6").Value = "1"
$ws.Range("F6").Value = "斤"
$ws.Range("G6").Value = "15"
$ws.Range("H6").Value = "15"
$ws.Range("I6").Value = "采购人"
$ws.Range("J6").Value = "收验货人"
$ws.Range("K6").Value = "供货人"

# ---- Row 7 ----
$ws.Range("A7").Value = "Tue Jul 10 2018 23:10:06 GMT+0800 (China Standard Time)"
$ws.Range("B7").Value = "瘦肉"
$ws.Range("E7").Value = "1"
$ws.Range("F7").Value = "斤"
$ws.Range("G7").Value = "15"
$ws.Range("H7").Value = "15"
$ws.Range("I7").Value = "采购人"
$ws.Range("J7").Value = "收验货人"
$ws.Range("K7").Value = "供货人"

# ---- Row 8 ----
$ws.Range("A8").Value = "Tue Jul 10 2018 23:13:14 GMT+0800 (China Standard Time)"
$ws.Range("B8").Value = "瘦肉"
$ws.Range("E8").Value = "1"
$ws.Range("F8").Value = "斤"
$ws.Range("G8").Value = "15"
$ws.Range("H8").Value = "15"
$ws.Range("I8").Value = "采购人"
$ws.Range("J8").Value = "收验货人"
$ws.Range("K8").Value = "供货人"

# ---- Row 9 ----
$ws.Range("A9").Value = "Tue Jul 10 2018 23:13:31 GMT+0800 (China Standard Time)"
$ws.Range("B9").Value = "瘦肉"
$ws.Range("C9").Value = "Tue Jul 10 2018 23:45:44 GMT+0800 (China Standard Time)"
$ws.Range("D9").Value = "Tue Jul 10 2018 23:45:46 GMT+0800 (China Standard Time)"
$ws.Range("E9").Value = "1"
$ws.Range("F9").Value = "斤"
$ws.Range("G9").Value = "15"
$ws.Range("H9").Value = "15"
$ws.Range("I9").Value = "采购人"
$ws.Range("J9").Value = "收验货人"
$ws.Range("K9").Value = "供货人"

# ---- Row 10 ----
$ws.Range("A10").Value = "Tue Jul 10 2018 23:45:36 GMT+0800 (China Standard Time)"
$ws.Range("B10").Value = "瘦肉"
$ws.Range("E10").Value = "1"
$ws.Range("F10").Value = "斤"
$ws.Range("G10").Value = "15"
$ws.Range("H10").Value = "15"
$ws.Range("I10").Value = "采购人"
$ws.Range("J10").Value = "收验货人"
$ws.Range("K10").Value = "供货人"
